$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.770.57'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.14%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.005.25'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.17%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '382.24'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.86%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '106.92'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.35%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.547'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.09%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.602'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.01%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.77'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.25%  '

$ws.Range('E11').Value = '  +0.56%  '

$ws.Range('E12').Value = '  +1.77%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.483.90'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.24%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.55'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.86%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.006.03'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.99%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.973'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.04%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.866.58'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.41%  '

$ws.Range('E19').Value = '  +2.31%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.46'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.14%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.17'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.28%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.75%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.05'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.25%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '264.70'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.93%  '

$ws.Range('E25').Value = '  +4.89%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.172'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.47%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.25'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +17.84%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.48'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.27%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '26.24'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.42%  '

$ws.Range('E30').Value = '  -0.10%  '

$ws.Range('E31').Value = '  +2.15%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.98'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.48%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '35.03'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.09%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '51.52'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.58%  '

$ws.Range('E35').Value = '  -2.87%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0455'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.88%  '

$ws.Range('E38').Value = '  +0.48%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '17.63'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.05%  '

$ws.Range('E40').Value = '  -5.25%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.87'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.57%  '

$ws.Range('E42').Value = '  +2.70%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '124.52'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.54%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '22.42'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.71%  '

$ws.Range('E45').Value = '  -1.11%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.278'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +17.48%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.062.77'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.11%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.37'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.72%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.30'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.55%  '

$ws.Range('E50').Value = '  +15.00%  '

$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.23'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.26%  '
